# Postdoc system database structure - "work on pdinit and system db structure"
#
# This script mirrors the manual edits made in Excel:
#  - storage sheet: leave a new cursor position (A13) behind
#  - archive sheet: no content change (just loses "last active tab" status)
#  - resources sheet:
#      * fix the second table's name from the duplicated "pgusers" to "pggroups"
#      * drop the "active" column from the "archives" table
#      * replace the "active" column in the "datasets" table with a new
#        "archive" column (FK into resources.archives)
#      * drop the "comma delimited" comment from the "applications" table's
#        features column
#      * rename "create" -> "createnew" in the "accesstokens" table
#      * drop the "comma delimited" comment from the accesstokens "features"
#        column, and reword the "content" column's comment
#      * end up with resources as the active/selected sheet, scrolled back
#        to the top, with C45 selected

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. storage sheet - just leave the selection on A13
# ---------------------------------------------------------------------------
$wsStorage = $wb.Worksheets.Item("storage")
$wsStorage.Activate()
$wsStorage.Range("A13").Select()

# ---------------------------------------------------------------------------
# 2. resources sheet - the actual schema edits
# ---------------------------------------------------------------------------
$wsRes = $wb.Worksheets.Item("resources")
$wsRes.Activate()

# second table header: "pgusers" (duplicate) -> "pggroups"
$wsRes.Range("A9").Value2 = "pggroups"

# "archives" table: remove the "active" field row entirely (was row 36)
$wsRes.Rows.Item(36).Delete()

# "applications" table: drop the "comma delimited" comment on "features"
# (now shifted up to row 52)
$wsRes.Range("D52").ClearContents()

# "accesstokens" table (now shifted up to rows 55-68):
#   reword the "content" column's comment
$wsRes.Range("D67").Value2 = "ALSO applies to content"

# drop the "comma delimited" comment on "features"
$wsRes.Range("D66").ClearContents()

#   "create" field renamed to "createnew"
$wsRes.Range("B61").Value2 = "createnew"

# "datasets" table: the old "active" row (now shifted up to row 44) becomes
# a new "archive" FK field
$wsRes.Range("B44").Value2 = "archive"
$wsRes.Range("C44").Value2 = "TEXT REFERENCES resources.archives(name)"

# leave the view scrolled to the top with C45 selected, as in the diff
$wsRes.Application.Goto($wsRes.Range("A2"))
$wsRes.Range("C45").Select()
